# Update countries & provincias Spain
#
# The "Pais" (COVID-19 dashboard) table is sorted by total cases
# (column B) descending. This refresh updates a handful of countries'
# figures; two of them (Polonia / Lituania) grew enough to overtake
# their neighbours in the ranking, so those rows shift down by one and
# the updated country is inserted just above them. All other changed
# rows are simple in-place value refreshes. The footer timestamp is
# also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (row 1) ---------------------------------------
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 10 de Octubre de 2020 a las 10:57"

# --- Rusia (row 7) ----------------------------------------------------
$ws.Cells.Item(7, 2).Value2 = 1285084
$ws.Cells.Item(7, 3).Value2 = 12846
$ws.Cells.Item(7, 4).Value2 = 1016202
$ws.Cells.Item(7, 5).Value2 = 246428
$ws.Cells.Item(7, 6).Value2 = 0
$ws.Cells.Item(7, 7).Value2 = 197
$ws.Cells.Item(7, 8).Value2 = 22454

# --- Filipinas (row 22) -----------------------------------------------
$ws.Cells.Item(22, 2).Value2 = 336926
$ws.Cells.Item(22, 3).Value2 = 2249
$ws.Cells.Item(22, 4).Value2 = 276094
$ws.Cells.Item(22, 5).Value2 = 54594
$ws.Cells.Item(22, 6).Value2 = 0
$ws.Cells.Item(22, 7).Value2 = 87
$ws.Cells.Item(22, 8).Value2 = 6238

# --- Indonesia (row 24) ------------------------------------------------
$ws.Cells.Item(24, 2).Value2 = 328952
$ws.Cells.Item(24, 3).Value2 = 4294
$ws.Cells.Item(24, 4).Value2 = 251481
$ws.Cells.Item(24, 5).Value2 = 65706
$ws.Cells.Item(24, 6).Value2 = 0
$ws.Cells.Item(24, 7).Value2 = 88
$ws.Cells.Item(24, 8).Value2 = 11765

# --- Polonia jumps above Panama / Republica Dominicana (rows 37-39) ---
$ws.Cells.Item(37, 1).Value2 = "Polonia"
$ws.Cells.Item(37, 2).Value2 = 121638
$ws.Cells.Item(37, 3).Value2 = 5300
$ws.Cells.Item(37, 4).Value2 = 78982
$ws.Cells.Item(37, 5).Value2 = 39684
$ws.Cells.Item(37, 6).Value2 = 0
$ws.Cells.Item(37, 7).Value2 = 53
$ws.Cells.Item(37, 8).Value2 = 2972

$ws.Cells.Item(38, 1).Value2 = "Panama"
$ws.Cells.Item(38, 2).Value2 = 118841
$ws.Cells.Item(38, 3).Value2 = 0
$ws.Cells.Item(38, 4).Value2 = 94962
$ws.Cells.Item(38, 5).Value2 = 21405
$ws.Cells.Item(38, 6).Value2 = 0
$ws.Cells.Item(38, 7).Value2 = 0
$ws.Cells.Item(38, 8).Value2 = 2474

$ws.Cells.Item(39, 1).Value2 = "Republica Dominicana"
$ws.Cells.Item(39, 2).Value2 = 117457
$ws.Cells.Item(39, 3).Value2 = 0
$ws.Cells.Item(39, 4).Value2 = 93061
$ws.Cells.Item(39, 5).Value2 = 22231
$ws.Cells.Item(39, 6).Value2 = 0
$ws.Cells.Item(39, 7).Value2 = 0
$ws.Cells.Item(39, 8).Value2 = 2165

# --- Singapur (row 62) --------------------------------------------------
$ws.Cells.Item(62, 2).Value2 = 57866
$ws.Cells.Item(62, 3).Value2 = 7
$ws.Cells.Item(62, 4).Value2 = 57675
$ws.Cells.Item(62, 5).Value2 = 164
$ws.Cells.Item(62, 6).Value2 = 0
$ws.Cells.Item(62, 7).Value2 = 0
$ws.Cells.Item(62, 8).Value2 = 27

# --- Austria (row 64) ----------------------------------------------------
$ws.Cells.Item(64, 2).Value2 = 54423
$ws.Cells.Item(64, 3).Value2 = 1235
$ws.Cells.Item(64, 4).Value2 = 42829
$ws.Cells.Item(64, 5).Value2 = 10742
$ws.Cells.Item(64, 6).Value2 = 0
$ws.Cells.Item(64, 7).Value2 = 10
$ws.Cells.Item(64, 8).Value2 = 852

# --- Croacia (row 91) -----------------------------------------------------
$ws.Cells.Item(91, 2).Value2 = 19932
$ws.Cells.Item(91, 3).Value2 = 486
$ws.Cells.Item(91, 4).Value2 = 16953
$ws.Cells.Item(91, 5).Value2 = 2662
$ws.Cells.Item(91, 6).Value2 = 0
$ws.Cells.Item(91, 7).Value2 = 4
$ws.Cells.Item(91, 8).Value2 = 317

# --- Eslovaquia (row 92) ---------------------------------------------------
$ws.Cells.Item(92, 2).Value2 = 18797
$ws.Cells.Item(92, 3).Value2 = 1887
$ws.Cells.Item(92, 4).Value2 = 5553
$ws.Cells.Item(92, 5).Value2 = 13183
$ws.Cells.Item(92, 6).Value2 = 0
$ws.Cells.Item(92, 7).Value2 = 4
$ws.Cells.Item(92, 8).Value2 = 61

# --- Lituania jumps above Cuba / Malaui (rows 121-123) ---------------------
$ws.Cells.Item(121, 1).Value2 = "Lituania"
$ws.Cells.Item(121, 2).Value2 = 5963
$ws.Cells.Item(121, 3).Value2 = 205
$ws.Cells.Item(121, 4).Value2 = 2751
$ws.Cells.Item(121, 5).Value2 = 3109
$ws.Cells.Item(121, 6).Value2 = 0
$ws.Cells.Item(121, 7).Value2 = 0
$ws.Cells.Item(121, 8).Value2 = 103

$ws.Cells.Item(122, 1).Value2 = "Cuba"
$ws.Cells.Item(122, 2).Value2 = 5943
$ws.Cells.Item(122, 3).Value2 = 0
$ws.Cells.Item(122, 4).Value2 = 5398
$ws.Cells.Item(122, 5).Value2 = 422
$ws.Cells.Item(122, 6).Value2 = 0
$ws.Cells.Item(122, 7).Value2 = 0
$ws.Cells.Item(122, 8).Value2 = 123

$ws.Cells.Item(123, 1).Value2 = "Malaui"
$ws.Cells.Item(123, 2).Value2 = 5813
$ws.Cells.Item(123, 3).Value2 = 0
$ws.Cells.Item(123, 4).Value2 = 4631
$ws.Cells.Item(123, 5).Value2 = 1002
$ws.Cells.Item(123, 6).Value2 = 0
$ws.Cells.Item(123, 7).Value2 = 0
$ws.Cells.Item(123, 8).Value2 = 180

# --- Estonia (row 141) -------------------------------------------------------
$ws.Cells.Item(141, 2).Value2 = 3846
$ws.Cells.Item(141, 3).Value2 = 37
$ws.Cells.Item(141, 4).Value2 = 2946
$ws.Cells.Item(141, 5).Value2 = 832
$ws.Cells.Item(141, 6).Value2 = 0
$ws.Cells.Item(141, 7).Value2 = 0
$ws.Cells.Item(141, 8).Value2 = 68

# --- Letonia (row 152) --------------------------------------------------------
$ws.Cells.Item(152, 2).Value2 = 2596
$ws.Cells.Item(152, 3).Value2 = 89
$ws.Cells.Item(152, 4).Value2 = 1322
$ws.Cells.Item(152, 5).Value2 = 1234
$ws.Cells.Item(152, 6).Value2 = 0
$ws.Cells.Item(152, 7).Value2 = 0
$ws.Cells.Item(152, 8).Value2 = 40
